$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-17 Friday" "2025-01-18 Saturday"
Replace-Text "959×5=" "322×8="
Replace-Text "129×9=" "495×4="
Replace-Text "459×8=" "358×6="
Replace-Text "944×9=" "381×8="
Replace-Text "565×4=" "789×3="
Replace-Text "718×6=" "738×4="
Replace-Text "437×3=" "156×2="
Replace-Text "526×8=" "426×5="
Replace-Text "327×5=" "587×2="
Replace-Text "240×5=" "140×2="
Replace-Text "231×5=" "516×2="
Replace-Text "337×2=" "485×3="
Replace-Text "707×9=" "848×4="
Replace-Text "651×4=" "259×9="
Replace-Text "389×6=" "630×9="
Replace-Text "435×4=" "529×2="
Replace-Text "696×9=" "451×6="
Replace-Text "337×4=" "540×4="
Replace-Text "883×5=" "681×3="
Replace-Text "238×3=" "774×8="
Replace-Text "855×4=" "487×7="
Replace-Text "799×5=" "390×8="
Replace-Text "692×7=" "972×4="
Replace-Text "224×9=" "949×5="
Replace-Text "754×4=" "496×7="
